# Append: 2026-02-12 19:00 JST
# The scraper ran again and produced two brand-new job postings which get
# inserted near the top of the list (as new rows 3 and 5), pushing the
# previously-seen postings down. Every row's "取得日時" (fetched-at)
# timestamp column is refreshed to the new run's timestamp, including rows
# whose job-posting content did not change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Make room for the two newly scraped rows -------------------------
# Final row 3 and final row 5 are brand-new entries; inserting in this order
# (first at 3, then at 5) reproduces the exact before->after row mapping:
#   old 2 -> 2 (unchanged content, timestamp refreshed)
#   old 3 -> 4
#   old 4 -> 6
#   old 5 -> 7
#   old 6 -> 8
#   old 7 -> 9
#   old 8 -> 10
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()

# --- 2. Clear out every existing hyperlink --------------------------------
# Row inserts in this runtime do not renumber the <hyperlinks> ref list, so
# the safest route is to drop every hyperlink now and recreate them fresh
# (step 4) once all the row data is in its final place.
$ws.Range("F2").Hyperlinks.Delete()

# --- 3. Write the final cell values for rows 2-10 -------------------------
$ws.Range("A2").Value = "2026-02-12 19:00:22"
$ws.Range("B2").Value = "【AI活用】市場調査・競合分析自動化機能の構築依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5489981"
$ws.Range("G2").Value = 388
$ws.Range("H2").Value = "🔥AI,Ai ◆自動化"

$ws.Range("A3").Value = "2026-02-12 19:00:22"
$ws.Range("B3").Value = "【急募】ビジネス向けAIエージェント開発支援のパートナー募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5490828"
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

$ws.Range("A4").Value = "2026-02-12 19:00:22"
$ws.Range("B4").Value = "【急募】ノーコードで実現するLINE×AI恋愛体験サービスMVP開発"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5490408"
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = "🔥AI,Ai ◆開発"

$ws.Range("A5").Value = "2026-02-12 19:00:22"
$ws.Range("B5").Value = "【注目】生成AIを活用したDX戦略の策定支援"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5490638"
$ws.Range("G5").Value = 310
$ws.Range("H5").Value = "🔥AI,Ai"

$ws.Range("A6").Value = "2026-02-12 19:00:22"
$ws.Range("B6").Value = "自動化システム"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5477084"
$ws.Range("G6").Value = 110
$ws.Range("H6").Value = "◆自動化"

$ws.Range("A7").Value = "2026-02-12 19:00:22"
$ws.Range("B7").Value = "【急募】FileMakerシステムのデバッグとレイアウト修正依頼"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5490478"
$ws.Range("G7").Value = 28
$ws.Range("H7").Value = ""

$ws.Range("A8").Value = "2026-02-12 19:00:22"
$ws.Range("B8").Value = "【長期】寝具ブランドのAmazon・楽天市場 運用代行パートナー募集"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5486471"
$ws.Range("G8").Value = 25
$ws.Range("H8").Value = ""

$ws.Range("A9").Value = "2026-02-12 19:00:22"
$ws.Range("B9").Value = "プロジェクトマネジメント"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5490062"
$ws.Range("G9").Value = 25
$ws.Range("H9").Value = ""

$ws.Range("A10").Value = "2026-02-12 19:00:22"
$ws.Range("B10").Value = "【SES経営者向け】事業立ち上げについてお話をお伺いできる方を募集します!"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5490407"
$ws.Range("G10").Value = 10
$ws.Range("H10").Value = ""

# --- 4. Recreate hyperlinks for the URL column, F2:F10 --------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5489981")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5490828")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5490408")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5490638")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5477084")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5490478")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5486471")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5490062")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5490407")

# Re-applying the named "Hyperlink" cell style normalises every F-column
# cell back onto the workbook's single shared Hyperlink style entry instead
# of leaving behind the ad-hoc duplicate style that Hyperlinks.Add() makes.
$ws.Range("F2:F10").Style = "Hyperlink"
